$wb = $excel.ActiveWorkbook

# --- 1. Employee sheet: add a "Skillset" column (C) ---
$wsEmployee = $wb.Worksheets.Item("Employee")

$wsEmployee.Range("C1").Font.Bold = $true
$wsEmployee.Range("C1").Value = "Skillset"

$employeeSkillset = @(1, 1, 2, 2, 3, 2, 3, 3, 2, 2, 3)
for ($i = 0; $i -lt $employeeSkillset.Length; $i++) {
    $wsEmployee.Cells.Item($i + 2, 3).Value = $employeeSkillset[$i]
}

# --- 2. Events sheet: insert two new columns (D, E) for Skillset1 / Skillset2 ---
$wsEvents = $wb.Worksheets.Item("Events")

$wsEvents.Columns("D:E").Insert()

$wsEvents.Range("D1").Value = "Skillset1"
$wsEvents.Range("E1").Value = "Skillset2"

$skillset1 = @(1, 0, 0, 1, 0, 0)
$skillset2 = @(4, 0, 0, 2, 0, 0)
for ($i = 0; $i -lt $skillset1.Length; $i++) {
    $row = $i + 2
    $wsEvents.Cells.Item($row, 4).Value = $skillset1[$i]
    $wsEvents.Cells.Item($row, 5).Value = $skillset2[$i]
}

# --- 3. View / selection bookkeeping (cosmetic, matches author's session) ---
$wsEmployee.Activate()
$wsEmployee.Range("B24").Select() | Out-Null

$wsEvents.Activate()
$wsEvents.Range("E14").Select() | Out-Null
$excel.ActiveWindow.Zoom = 75
